# Auto-generated Excel COM-interop script
# Applies the cryptos list update as described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.937.97"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "1.822.84"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D4").Value = "'0.9918"
$ws.Range("E4").Value = "  -0.75%  "
$ws.Range("D5").Value = "'243.42"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").Value = "'0.6294"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "'0.9967"
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("D8").Value = "'0.07460"
$ws.Range("E8").Value = "  -1.89%  "
$ws.Range("D9").Value = "'0.2932"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").Value = "'23.01"
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("D11").Value = "'0.07678"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "1.826.63"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").Value = "'4.978"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").Value = "'0.6662"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("D15").Value = "'82.88"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").Value = "'0.000009585"
$ws.Range("E16").Value = "  +1.90%  "
$ws.Range("D17").Value = "'6.036"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").Value = "28.989.98"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").Value = "'12.54"
$ws.Range("E19").Value = "  +1.71%  "
$ws.Range("D20").Value = "'225.47"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "'0.9958"
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("D22").Value = "'7.116"
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("D23").Value = "'0.9943"
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("D24").Value = "'160.07"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").Value = "'0.1408"
$ws.Range("E25").Value = "  +3.30%  "
$ws.Range("D26").Value = "'8.480"
$ws.Range("E26").Value = "  +0.69%  "
$ws.Range("D27").Value = "'17.89"
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("D28").Value = "'1.498"
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").Value = "'4.116"
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("D30").Value = "'4.049"
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("D31").Value = "'0.05434"
$ws.Range("E31").Value = "  +4.37%  "
$ws.Range("D32").Value = "'1.196"
$ws.Range("D33").Value = "'1.850"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("D34").Value = "'0.7409"
$ws.Range("E34").Value = "  +1.40%  "
$ws.Range("D35").Value = "'1.133"
$ws.Range("E35").Value = "  -1.52%  "
$ws.Range("D36").Value = "'2.614"
$ws.Range("E36").Value = "  +0.63%  "
$ws.Range("D37").Value = "1.240.43"
$ws.Range("E37").Value = "  -2.52%  "
$ws.Range("D38").Value = "'2.743"
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("D39").Value = "'0.01776"
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("D40").Value = "'6.645"
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("D41").Value = "'0.8986"
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("D42").Value = "'0.9972"
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("D43").Value = "'101.23"
$ws.Range("E43").Value = "  -0.47%  "
$ws.Range("D44").Value = "1.973.87"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'64.92"
$ws.Range("E45").Value = "  +1.39%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.00000000122"
$ws.Range("E46").Value = "  +1.87%  "
$ws.Range("D47").Value = "'0.5065"
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("D48").Value = "'0.4036"
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("D49").Value = "'8.949"
$ws.Range("E49").Value = "  +1.25%  "
$ws.Range("D50").Value = "'1.654"
$ws.Range("E50").Value = "  +1.28%  "
$ws.Range("D51").Value = "'0.07163"
$ws.Range("E51").Value = "  -0.44%  "
